# [export] Display HoursWorked for exported warrants
#
# The "Detail" sheet gets a new trailing column (J) with header
# "Hours Worked" so the exported warrant detail shows hours worked
# alongside the existing payroll figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Detail")

# Copy the format/style of the last existing header cell (I1, "Net Pay")
# onto the new header cell J1, then give it its own text.
$ws.Range("I1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "Hours Worked"

# Size the new column similarly to the other bestFit header columns.
$ws.Columns.Item(10).ColumnWidth = 15.1666666666667
